# Applies the latest scheduled-runner price/profit recalculation to the
# Leve profit tables on each crafting-class worksheet (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Only numeric market-price/profit columns (H-N) for
# specific leve rows are touched; all other data is left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 376.33334
$ws.Range("I4").Value = 376.33334
$ws.Range("K4").Value = 376.33334
$ws.Range("M4").Value = -262.33334

# Row 21
$ws.Range("H21").Value = 17
$ws.Range("I21").Value = 17
$ws.Range("K21").Value = 17
$ws.Range("M21").Value = 451

# Row 23
$ws.Range("H23").Value = 17
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = 17
$ws.Range("M23").Value = 217

# Row 70
$ws.Range("H70").Value = 2937.25
$ws.Range("I70").Value = 2899.8
$ws.Range("K70").Value = 8699.400000000001
$ws.Range("M70").Value = -8429.400000000001

# Row 73
$ws.Range("H73").Value = 2937.25
$ws.Range("I73").Value = 2899.8
$ws.Range("K73").Value = 8699.400000000001
$ws.Range("M73").Value = -7763.400000000001

# Row 80
$ws.Range("H80").Value = 502
$ws.Range("I80").Value = 502
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1506
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -508
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 502
$ws.Range("I83").Value = 502
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 4518
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 474
$ws.Range("N83").ClearContents()

# Row 107
$ws.Range("H107").Value = 666.1
$ws.Range("I107").Value = 684.8333
$ws.Range("J107").Value = 497.5
$ws.Range("K107").Value = 684.8333
$ws.Range("L107").Value = 497.5
$ws.Range("M107").Value = 1235.1667
$ws.Range("N107").Value = -4337.5

# Row 137
$ws.Range("H137").Value = 1056.5264
$ws.Range("I137").Value = 840.1875
$ws.Range("K137").Value = 2520.5625
$ws.Range("M137").Value = 29.4375

$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 6049
$ws.Range("I25").Value = 2600
$ws.Range("J25").Value = 9498
$ws.Range("K25").Value = 2600
$ws.Range("L25").Value = 9498
$ws.Range("M25").Value = -2198
$ws.Range("N25").Value = -10302

# Row 74
$ws.Range("H74").Value = 1583
$ws.Range("I74").Value = 1583
$ws.Range("K74").Value = 1583
$ws.Range("M74").Value = -709

# Row 77
$ws.Range("H77").Value = 1583
$ws.Range("I77").Value = 1583
$ws.Range("K77").Value = 7915
$ws.Range("M77").Value = -3547

# Row 102
$ws.Range("H102").Value = 57191.11
$ws.Range("I102").Value = 67962.664
$ws.Range("K102").Value = 67962.664
$ws.Range("M102").Value = -66340.664

# Row 122
$ws.Range("H122").Value = 1879.6
$ws.Range("I122").Value = 1849.75
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 5549.25
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -3099.25
$ws.Range("N122").Value = -10897

# Row 132
$ws.Range("H132").Value = 1893.1538
$ws.Range("I132").Value = 1691.909
$ws.Range("K132").Value = 5075.727000000001
$ws.Range("M132").Value = -2545.727000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 3421.52
$ws.Range("I99").Value = 3403.8262
$ws.Range("K99").Value = 3403.8262
$ws.Range("M99").Value = -1905.8262

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3019.611
$ws.Range("I31").Value = 1676.0714
$ws.Range("J31").Value = 7722
$ws.Range("K31").Value = 1676.0714
$ws.Range("L31").Value = 7722
$ws.Range("M31").Value = -1381.0714
$ws.Range("N31").Value = -8312

# Row 34
$ws.Range("H34").Value = 3019.611
$ws.Range("I34").Value = 1676.0714
$ws.Range("J34").Value = 7722
$ws.Range("K34").Value = 1676.0714
$ws.Range("L34").Value = 7722
$ws.Range("M34").Value = -1474.0714
$ws.Range("N34").Value = -8126

# Row 99
$ws.Range("H99").Value = 2526250
$ws.Range("I99").Value = 2550000
$ws.Range("J99").Value = 2502500
$ws.Range("K99").Value = 2550000
$ws.Range("L99").Value = 2502500
$ws.Range("M99").Value = -2548502
$ws.Range("N99").Value = -2505496

# Row 126
$ws.Range("H126").Value = 2526250
$ws.Range("I126").Value = 2550000
$ws.Range("J126").Value = 2502500
$ws.Range("K126").Value = 7650000
$ws.Range("L126").Value = 7507500
$ws.Range("M126").Value = -7647530
$ws.Range("N126").Value = -7512440

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 2001.8889
$ws.Range("J22").Value = 2001.8889
$ws.Range("L22").Value = 6005.6667
$ws.Range("N22").Value = -6343.6667

# Row 27
$ws.Range("H27").Value = 2001.8889
$ws.Range("J27").Value = 2001.8889
$ws.Range("L27").Value = 6005.6667
$ws.Range("N27").Value = -6209.6667

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6000
$ws.Range("J70").Value = 6000
$ws.Range("L70").Value = 6000
$ws.Range("N70").Value = -6540

# Row 73
$ws.Range("H73").Value = 6000
$ws.Range("J73").Value = 6000
$ws.Range("L73").Value = 6000
$ws.Range("N73").Value = -7872

# Row 117
$ws.Range("H117").Value = 23750
$ws.Range("J117").Value = 23750
$ws.Range("L117").Value = 23750
$ws.Range("N117").Value = -30634

# Row 122
$ws.Range("H122").Value = 7406.273
$ws.Range("I122").Value = 6624
$ws.Range("K122").Value = 19872
$ws.Range("M122").Value = -17422

# Row 132
$ws.Range("H132").Value = 1484.4
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4946.7
$ws.Range("I61").Value = 4433.5
$ws.Range("K61").Value = 4433.5
$ws.Range("M61").Value = -4231.5

# Row 93
$ws.Range("H93").Value = 1962.5
$ws.Range("I93").Value = 2283.3333
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 2283.3333
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -1035.3333
$ws.Range("N93").Value = -3496

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 113
$ws.Range("H113").Value = 4946.7
$ws.Range("I113").Value = 4433.5
$ws.Range("K113").Value = 4433.5
$ws.Range("M113").Value = -2263.5

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Row 41
$ws.Range("H41").Value = 34999.668
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20780

# Row 136
$ws.Range("H136").Value = 3249.6667
$ws.Range("I136").Value = 3249.6667
$ws.Range("K136").Value = 9749.000100000001
$ws.Range("M136").Value = -7199.000100000001
